$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.791.80"
$ws.Range("E2").Value = "  -3.95%  "
$ws.Range("D3").Value = "1.817.00"
$ws.Range("E3").Value = "  -2.95%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "277.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -7.83%  "
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("E7").Value = "  -4.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3514"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -6.23%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.56"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.12%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06661"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -7.17%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.97"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -7.21%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.8322"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -6.19%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07899"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.21%  "
$ws.Range("D14").Value = "1.818.51"
$ws.Range("E14").Value = "  -3.14%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.078"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.50%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "87.64"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -6.03%  "
$ws.Range("E17").Value = "  -0.01%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.12"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000008020"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -6.02%  "
$ws.Range("E20").Value = "  -0.07%  "
$ws.Range("D21").Value = "25.847.99"
$ws.Range("E21").Value = "  -3.90%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.730"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.84%  "
$ws.Range("E23").Value = "  -6.29%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.087"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.82%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "141.84"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.55%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.182"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.36%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.669"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.08"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.23%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "109.50"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.16%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.346"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -8.25%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.235"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -7.65%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08837"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.14%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04867"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.46%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7300"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -8.65%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.133"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.75%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.878"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.79%  "
$ws.Range("E37").Value = "  -0.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.9996"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.05%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.338"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -9.74%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5231"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -13.33%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01850"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.25%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9551"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -10.94%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.207"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.20%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "111.21"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.99%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.058"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -9.10%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9999"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4595"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -10.89%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1363"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -8.99%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "36.72"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.47%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.257"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -7.07%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.502"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -7.52%  "
